# DOMA-746 add executor & assignee excel export mappers
#
# The ticket-analytics export template referenced a single `{d.ticket[...]}`
# collection; the exporter now feeds a `tickets` collection instead, so every
# `{d.ticket[i]...}` / `{d.ticket[i + 1]...}` placeholder on the sheet is
# renamed to `{d.tickets[i]...}` / `{d.tickets[i + 1]...}`.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 — the "i"-th ticket row of placeholders.
$ws.Range("A2").Value = "{d.tickets[i].address}"
$ws.Range("B2").Value = "{d.tickets[i].processing}"
$ws.Range("C2").Value = "{d.tickets[i].completed}"
$ws.Range("D2").Value = "{d.tickets[i].canceled}"
$ws.Range("E2").Value = "{d.tickets[i].deferred}"
$ws.Range("F2").Value = "{d.tickets[i].closed}"
$ws.Range("G2").Value = "{d.tickets[i].new_or_reopened}"

# Row 3 — the "i + 1"-th ticket row of placeholders.
$ws.Range("A3").Value = "{d.tickets[i + 1].address}"
$ws.Range("B3").Value = "{d.tickets[i + 1].processing}"
$ws.Range("C3").Value = "{d.tickets[i + 1].completed}"
$ws.Range("D3").Value = "{d.tickets[i + 1].canceled}"
$ws.Range("E3").Value = "{d.tickets[i + 1].deferred}"
$ws.Range("F3").Value = "{d.tickets[i + 1].closed}"
$ws.Range("G3").Value = "{d.tickets[i + 1].new_or_reopened}"

# The editor's cursor ended up parked on D28 when the file was last saved.
$ws.Range("D28").Select()
